$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp note in A1
$ws.Range("A1").Value = "Datos actualizados a 10 de Octubre de 2020 a las 13:31"

# Row 4
$ws.Range("B4").Value = 7895498
$ws.Range("C4").Value = 1020
$ws.Range("D4").Value = 5065054
$ws.Range("E4").Value = 2611775
$ws.Range("G4").Value = 21
$ws.Range("H4").Value = 218669

# Row 16
$ws.Range("B16").Value = 496253
$ws.Range("C16").Value = 3875
$ws.Range("D16").Value = 403950
$ws.Range("E16").Value = 64010
$ws.Range("G16").Value = 195
$ws.Range("H16").Value = 28293

# Row 40
$ws.Range("B40").Value = 110568
$ws.Range("C40").Value = 492
$ws.Range("D40").Value = 102722
$ws.Range("E40").Value = 7191
$ws.Range("G40").Value = 6
$ws.Range("H40").Value = 655

# Row 43
$ws.Range("A43").Value = "Nepal"
$ws.Range("B43").Value = 105684
$ws.Range("C43").Value = 5008
$ws.Range("D43").Value = 74252
$ws.Range("E43").Value = 30818
$ws.Range("G43").Value = 14
$ws.Range("H43").Value = 614

# Row 44
$ws.Range("A44").Value = "Emiratos Arabes Unidos"
$ws.Range("B44").Value = 105133
$ws.Range("C44").Value = 1129
$ws.Range("D44").Value = 95973
$ws.Range("E44").Value = 8717
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 443

# Row 45
$ws.Range("A45").Value = "Egipto"
$ws.Range("B45").Value = 104262
$ws.Range("D45").Value = 97592
$ws.Range("E45").Value = 641
$ws.Range("H45").Value = 6029

# Row 46
$ws.Range("A46").Value = "Oman"
$ws.Range("B46").Value = 104129
$ws.Range("D46").Value = 91731
$ws.Range("E46").Value = 11389
$ws.Range("H46").Value = 1009

# Row 59
$ws.Range("A59").Value = "Uzbekistan"
$ws.Range("B59").Value = 60562
$ws.Range("C59").Value = 220
$ws.Range("D59").Value = 57330
$ws.Range("E59").Value = 2732
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 500

# Row 60
$ws.Range("A60").Value = "Suiza"
$ws.Range("B60").Value = 60368
$ws.Range("D60").Value = 48400
$ws.Range("E60").Value = 9880
$ws.Range("H60").Value = 2088

# Row 93
$ws.Range("B93").Value = 16702
$ws.Range("C93").Value = 26
$ws.Range("D93").Value = 16012
$ws.Range("E93").Value = 453

# Row 97
$ws.Range("A97").Value = "Malasia"
$ws.Range("B97").Value = 15096
$ws.Range("C97").Value = 374
$ws.Range("D97").Value = 10780
$ws.Range("E97").Value = 4161
$ws.Range("G97").Value = 3
$ws.Range("H97").Value = 155

# Row 98
$ws.Range("A98").Value = "Albania"
$ws.Range("B98").Value = 15066
$ws.Range("D98").Value = 9304
$ws.Range("E98").Value = 5349
$ws.Range("H98").Value = 413

# Row 114
$ws.Range("A114").Value = "Eslovenia"
$ws.Range("B114").Value = 8252
$ws.Range("C114").Value = 380
$ws.Range("D114").Value = 5024
$ws.Range("E114").Value = 3061
$ws.Range("G114").Value = 2
$ws.Range("H114").Value = 167

# Row 115
$ws.Range("A115").Value = "Zimbabue"
$ws.Range("B115").Value = 7994
$ws.Range("D115").Value = 6474
$ws.Range("E115").Value = 1291
$ws.Range("H115").Value = 229

# Row 142
$ws.Range("A142").Value = "Malta"
$ws.Range("B142").Value = 3681
$ws.Range("C142").Value = 100
$ws.Range("D142").Value = 2937
$ws.Range("E142").Value = 703
$ws.Range("H142").Value = 41

# Row 143
$ws.Range("A143").Value = "Tailandia"
$ws.Range("B143").Value = 3634
$ws.Range("C143").Value = 6
$ws.Range("D143").Value = 3445
$ws.Range("E143").Value = 130
$ws.Range("H143").Value = 59

# Row 144
$ws.Range("A144").Value = "Gambia"
$ws.Range("B144").Value = 3621
$ws.Range("D144").Value = 2489
$ws.Range("E144").Value = 1015
$ws.Range("H144").Value = 117

# Row 181
$ws.Range("B181").Value = 476
$ws.Range("C181").Value = 8
$ws.Range("D181").Value = 405
$ws.Range("E181").Value = 71

# Row 196
$ws.Range("B196").Value = 142
$ws.Range("C196").Value = 3
$ws.Range("E196").Value = 24
